# Rician bias 보정 전체
# Apply corrected SNR values to each sheet's data table (rows 2-7, columns A-C)

$wb = $excel.ActiveWorkbook

$sheetData = @{
    "Complex_SNR" = @(
        @(1.886208295822144, 9.377218246459961, 7.491009712219238),
        @(-1.054936528205872, 6.329257488250732, 7.384193897247314),
        @(-4.264941215515137, 2.887811422348022, 7.152752876281738),
        @(-7.546406269073486, -0.6539157629013062, 6.892490386962891),
        @(-10.55385971069336, -4.003417015075684, 6.550442695617676),
        @(-13.02655410766602, -6.857790470123291, 6.168763637542725)
    )
    "Magnitude_SNR" = @(
        @(8.391603469848633, 15.32101249694824, 6.929409027099609),
        @(7.670703411102295, 14.7795581817627, 7.1088547706604),
        @(6.96019172668457, 13.95021629333496, 6.990024566650391),
        @(6.286073684692383, 13.01476860046387, 6.728694915771484),
        @(5.618660926818848, 11.93073844909668, 6.312077522277832),
        @(4.953777313232422, 10.77796459197998, 5.824187278747559)
    )
    "Corrected_SNR" = @(
        @(7.935194492340088, 15.01033878326416, 7.075144290924072),
        @(7.150395393371582, 14.36503887176514, 7.214643478393555),
        @(6.372131824493408, 13.40629577636719, 7.034163951873779),
        @(5.629256248474121, 12.3375883102417, 6.708332061767578),
        @(4.897980690002441, 11.15977764129639, 6.261796951293945),
        @(4.167629718780518, 9.93265438079834, 5.765024662017822)
    )
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $rowIndex = $i + 2
        $vals = $rows[$i]
        $ws.Cells.Item($rowIndex, 1).Value = $vals[0]
        $ws.Cells.Item($rowIndex, 2).Value = $vals[1]
        $ws.Cells.Item($rowIndex, 3).Value = $vals[2]
    }
}

$wb.Save()
